$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text formatting on price cells whose new value would otherwise be
# auto-detected as a number by Excel's input parser, so they stay text
# (matching the original inline-string cell content).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the updated coin / link / price / volume values.
$ws.Range("D2").Value = "42.977.99"
$ws.Range("E2").Value = "  +0.42%  "
$ws.Range("D3").Value = "2.366.36"
$ws.Range("E3").Value = "  +2.41%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "302.45"
$ws.Range("E5").Value = "  +0.35%  "
$ws.Range("D6").Value = "95.78"
$ws.Range("E6").Value = "  +0.63%  "
$ws.Range("E7").Value = "  -0.33%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "0.488"
$ws.Range("E9").Value = "  -0.70%  "
$ws.Range("D10").Value = "33.99"
$ws.Range("E10").Value = "  -0.39%  "
$ws.Range("D11").Value = "0.0785"
$ws.Range("E11").Value = "  +0.33%  "
$ws.Range("E12").Value = "  +3.19%  "
$ws.Range("D13").Value = "18.39"
$ws.Range("E13").Value = "  -3.12%  "
$ws.Range("D14").Value = "6.74"
$ws.Range("E14").Value = "  +0.09%  "
$ws.Range("D15").Value = "2.730.15"
$ws.Range("E15").Value = "  +2.15%  "
$ws.Range("D16").Value = "2.392.19"
$ws.Range("E16").Value = "  +2.84%  "
$ws.Range("D17").Value = "0.795"
$ws.Range("E17").Value = "  +0.95%  "
$ws.Range("D18").Value = "42.950.08"
$ws.Range("D19").Value = "11.93"
$ws.Range("E19").Value = "  -2.14%  "
$ws.Range("D20").Value = "6.25"
$ws.Range("E20").Value = "  +1.98%  "
$ws.Range("D21").Value = "0.0₃0885"
$ws.Range("E21").Value = "  -0.54%  "
$ws.Range("D22").Value = "67.96"
$ws.Range("E22").Value = "  +0.41%  "
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("E24").Value = "  -4.31%  "
$ws.Range("E25").Value = "  +0.04%  "
$ws.Range("E26").Value = "  +0.20%  "
$ws.Range("D27").Value = "24.61"
$ws.Range("E27").Value = "  +1.39%  "
$ws.Range("D28").Value = "2.37"
$ws.Range("E28").Value = "  +0.45%  "
$ws.Range("E29").Value = "  +1.56%  "
$ws.Range("D30").Value = "31.47"
$ws.Range("E30").Value = "  -1.88%  "
$ws.Range("D31").Value = "1.00"
$ws.Range("E31").Value = "  -0.03%  "
$ws.Range("D32").Value = "5.02"
$ws.Range("E32").Value = "  +0.40%  "
$ws.Range("D33").Value = "17.28"
$ws.Range("E33").Value = "  -2.15%  "
$ws.Range("D34").Value = "0.0721"
$ws.Range("E34").Value = "  +3.54%  "
$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D35").Value = "0.104"
$ws.Range("E35").Value = "  +4.13%  "
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").Value = "1.84"
$ws.Range("E36").Value = "  +3.72%  "
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D37").Value = "4.34"
$ws.Range("E37").Value = "  -2.20%  "
$ws.Range("B38").Value = "Monero"
$ws.Range("C38").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D38").Value = "125.45"
$ws.Range("E38").Value = "  -24.54%  "
$ws.Range("D39").Value = "2.28"
$ws.Range("E39").Value = "  -2.16%  "
$ws.Range("E40").Value = "  +2.64%  "
$ws.Range("E41").Value = "  -0.49%  "
$ws.Range("D42").Value = "21.35"
$ws.Range("E42").Value = "  +2.01%  "
$ws.Range("D43").Value = "1.936.87"
$ws.Range("E43").Value = "  +0.53%  "
$ws.Range("E44").Value = "  +0.25%  "
$ws.Range("E45").Value = "  +2.38%  "
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").Value = "9.18"
$ws.Range("E46").Value = "  -9.43%  "
$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").Value = "2.71"
$ws.Range("E47").Value = "  -0.78%  "
$ws.Range("D48").Value = "2.588.39"
$ws.Range("E48").Value = "  +1.90%  "
$ws.Range("E49").Value = "  +2.16%  "
$ws.Range("D50").Value = "1.15"
$ws.Range("E50").Value = "  +1.73%  "
$ws.Range("D51").Value = "71.41"
$ws.Range("E51").Value = "  -0.86%  "
